# Append new ActivityLog entries (rows 13-22) recorded after the previous
# save, restoring the local changes that existed before pulling updates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Timestamp, B=Admin, C=Action, D=Account Affected (text), E=Amount (number)
$rows = @(
    @("2025-10-22 19:59:45", "Sumit", "Deposit",  "12344",  10000),
    @("2025-10-22 19:59:55", "Sumit", "Deposit",  "12344",  10000000),
    @("2025-10-22 20:00:09", "Sumit", "Deposit",  "12344",  200000000),
    @("2025-10-22 20:00:21", "Sumit", "Withdraw", "123456", 1),
    @("2025-10-22 20:00:29", "Sumit", "Withdraw", "12344",  2000000),
    @("2025-10-22 20:00:37", "Sumit", "Withdraw", "12344",  20000000),
    @("2025-10-22 20:00:47", "Sumit", "Withdraw", "12344",  19000000),
    @("2025-10-22 20:01:13", "Sumit", "Withdraw", "12344",  179010425),
    @("2025-10-22 20:01:27", "Sumit", "Deposit",  "12344",  1000),
    @("2025-10-22 23:16:27", "Sumit", "Withdraw", "12344",  999)
)

$startRow = 13
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $entry = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
    # Leading apostrophe forces the numeric-looking account id to stay text,
    # matching the other "Account Affected" entries in the sheet.
    $ws.Cells.Item($r, 4).Value = "'" + $entry[3]
    $ws.Cells.Item($r, 5).Value = $entry[4]
}
